# Updates cryptos list figures (price & 1h volume change) and
# re-labels a few reordered coin rows, per the scraped data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '34.565.47'
$ws.Range('E2').Value = '  +0.10%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '1.818.23'
$ws.Range('E3').Value = '  +0.30%  '
$ws.Range('E4').Value = '  -0.04%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '226.58'
$ws.Range('E5').Value = '  -0.97%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '0.596'
$ws.Range('E6').Value = '  +2.65%  '
$ws.Range('E7').Value = '  -0.07%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '38.50'
$ws.Range('E8').Value = '  +6.74%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '0.292'
$ws.Range('E9').Value = '  -3.89%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '0.0681'
$ws.Range('E10').Value = '  -2.79%  '
$ws.Range('E11').Value = '  +1.66%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '2.081.93'
$ws.Range('E12').Value = '  +0.41%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '11.33'
$ws.Range('E13').Value = '  -0.65%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '1.815.50'
$ws.Range('E14').Value = '  +0.72%  '
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D15').Value = '0.639'
$ws.Range('E15').Value = '  -1.90%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '4.47'
$ws.Range('E16').Value = '  -2.20%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '34.573.50'
$ws.Range('E17').Value = '  +0.15%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '69.00'
$ws.Range('E18').Value = '  -0.71%  '
$ws.Range('D19').NumberFormat = "@"
$ws.Range('D19').Value = '244.93'
$ws.Range('E19').Value = '  -1.53%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '0.0₃0778'
$ws.Range('E20').Value = '  -2.92%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '11.29'
$ws.Range('E21').Value = '  -3.12%  '
$ws.Range('E22').Value = '  -0.12%  '
$ws.Range('E23').Value = '  -1.24%  '
$ws.Range('E24').Value = '  +4.86%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '171.43'
$ws.Range('E25').Value = '  -0.36%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '7.89'
$ws.Range('E26').Value = '  -1.87%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '17.84'
$ws.Range('E27').Value = '  +5.58%  '
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '0.121'
$ws.Range('E28').Value = '  +1.04%  '
$ws.Range('E29').Value = '  -0.04%  '
$ws.Range('B30').Value = 'Filecoin'
$ws.Range('C30').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D30').Value = '3.82'
$ws.Range('E30').Value = '  -1.98%  '
$ws.Range('B31').Value = 'PancakeSwap'
$ws.Range('C31').Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '1.24'
$ws.Range('E31').Value = '  -0.91%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '3.92'
$ws.Range('E32').Value = '  -4.56%  '
$ws.Range('B33').Value = 'Hedera'
$ws.Range('C33').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '0.0523'
$ws.Range('E33').Value = '  -2.32%  '
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D34').Value = '1.85'
$ws.Range('E34').Value = '  -0.60%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.367.48'
$ws.Range('E35').Value = '  -2.68%  '
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D36').Value = '0.655'
$ws.Range('E36').Value = '  -4.16%  '
$ws.Range('E37').Value = '  -0.09%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '2.38'
$ws.Range('E38').Value = '  -5.51%  '
$ws.Range('E39').Value = '  -1.77%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '0.962'
$ws.Range('E40').Value = '  -1.09%  '
$ws.Range('B41').Value = 'HuobiToken'
$ws.Range('C41').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '2.43'
$ws.Range('E41').Value = '  +1.25%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '81.94'
$ws.Range('E42').Value = '  -2.67%  '
$ws.Range('B43').Value = 'WEMIXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '1.21'
$ws.Range('E43').Value = '  +2.82%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '2.79'
$ws.Range('E44').Value = '  -1.31%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '13.95'
$ws.Range('E45').Value = '  +3.87%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0508'
$ws.Range('E46').Value = '  +0.27%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '1.981.72'
$ws.Range('E47').Value = '  +0.42%  '
$ws.Range('E48').Value = '  -3.64%  '
$ws.Range('E49').Value = '  -0.15%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '103.21'
$ws.Range('E50').Value = '  -2.60%  '
$ws.Range('B51').Value = 'BabyDogeCoin'
$ws.Range('C51').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D51').NumberFormat = "@"
$ws.Range('D51').Value = '0.0₆0122'
$ws.Range('E51').Value = '  -6.31%  '
